$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 corresponds to file_name = metrics_sim_with_priors.json
# Correcting relevance markers: update recall/tnr/confusion-matrix metrics

$ws.Range("C3").Value = 0.9736842105263158
$ws.Range("D3").Value = 1

$ws.Range("H3").Value = 0.8928257042253521
$ws.Range("I3").Value = 0.02152685308477586
$ws.Range("J3").Value = 0.868421052631579
$ws.Range("K3").Value = 116.5

$ws.Range("Q3").Value = 33
$ws.Range("R3").Value = 61
$ws.Range("S3").Value = 76
$ws.Range("T3").Value = 93
$ws.Range("U3").Value = 121

$ws.Range("V3").Value = 4473
$ws.Range("W3").Value = 4445
$ws.Range("X3").Value = 4430
$ws.Range("Y3").Value = 4413
$ws.Range("Z3").Value = 4385

$ws.Range("AF3").Value = 0.992676
$ws.Range("AG3").Value = 0.986462
$ws.Range("AH3").Value = 0.983134
$ws.Range("AI3").Value = 0.979361
$ws.Range("AJ3").Value = 0.973147
